$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shifted/new data for rows 181-223 (columns D, J, K, L, M, N, O, P, Q)
$rows = @(
    @{Row=181; D=44637; J=80; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=182; D=44208; J=100; K=5000; L=5000; M=5000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2500; Q=2},
    @{Row=183; D=44355; J=160; K=3500; L=3500; M=3500; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1167; Q=3},
    @{Row=184; D=44530; J=120; K=6000; L=6000; M=6000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=3000; Q=2},
    @{Row=185; D=44530; J=180; K=5500; L=5500; M=5500; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1833; Q=3},
    @{Row=186; D=44483; J=80; K=4500; L=4500; M=4500; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1500; Q=3},
    @{Row=187; D=44294; J=70; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=188; D=44617; J=180; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=189; D=44264; J=80; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=190; D=44264; J=80; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=191; D=44232; J=100; K=5000; L=5000; M=5000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2500; Q=2},
    @{Row=192; D=44279; J=30; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=193; D=44330; J=180; K=4000; L=4000; M=4000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1333; Q=3},
    @{Row=194; D=44504; J=60; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=195; D=44257; J=100; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=196; D=44257; J=100; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=197; D=44301; J=80; K=4500; L=4500; M=4500; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1500; Q=3},
    @{Row=198; D=44370; J=20; K=4000; L=4000; M=4000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1333; Q=3},
    @{Row=199; D=44487; J=90; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=200; D=44385; J=30; K=4000; L=4000; M=4000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1333; Q=3},
    @{Row=201; D=44236; J=100; K=5000; L=5000; M=5000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2500; Q=2},
    @{Row=202; D=44221; J=100; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=203; D=44272; J=20; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=204; D=44229; J=200; K=4000; L=5000; M=4500; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2250; Q=2},
    @{Row=205; D=44214; J=50; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=206; D=44299; J=80; K=5000; L=5000; M=5000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2500; Q=2},
    @{Row=207; D=44299; J=140; K=4500; L=4500; M=4500; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1500; Q=3},
    @{Row=208; D=44610; J=150; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=209; D=44312; J=20; K=5000; L=5000; M=5000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2500; Q=2},
    @{Row=210; D=44522; J=60; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=211; D=44277; J=80; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=212; D=44258; J=20; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=213; D=44390; J=180; K=4500; L=4500; M=4500; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1500; Q=3},
    @{Row=214; D=44349; J=30; K=3500; L=3500; M=3500; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1167; Q=3},
    @{Row=215; D=44285; J=100; K=4000; L=4000; M=4000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=2000; Q=2},
    @{Row=216; D=44285; J=150; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=217; D=44498; J=160; K=4500; L=4500; M=4500; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1500; Q=3},
    @{Row=218; D=44418; J=180; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=219; D=44595; J=80; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3},
    @{Row=220; D=44628; J=180; K=6000; L=6000; M=6000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=2000; Q=3},
    @{Row=221; D=44552; J=20; K=7000; L=7000; M=7000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=3500; Q=2},
    @{Row=222; D=44544; J=80; K=6000; L=6000; M=6000; N="`$/docena de atados (2 kilos)"; O="Región de La Araucanía"; P=3000; Q=2},
    @{Row=223; D=44544; J=160; K=5000; L=5000; M=5000; N="`$/docena de atados (3 kilos)"; O="Región Metropolitana"; P=1667; Q=3}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # N: Unidad de comercializacion
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio $/Kg
    $ws.Cells.Item($r.Row, 17).Value = $r.Q   # Q: Kg o Unidades
}

# Row 223 is brand new: fill the constant columns that were identical across the whole block
$ws.Cells.Item(223, 1).Value = 4                                      # A: Mercado ID
$ws.Cells.Item(223, 2).Value = "Feria Lagunitas de Puerto Montt"      # B: Mercado
$ws.Cells.Item(223, 3).Value = "Los Lagos"                            # C: Region
$ws.Cells.Item(223, 5).Value = 10                                     # E: Codreg
$ws.Cells.Item(223, 6).Value = 100112044                             # F: Categoria ID
$ws.Cells.Item(223, 7).Value = "Perejil"                              # G: Categoria
$ws.Cells.Item(223, 8).Value = "Sin especificar"                      # H: Variedad
$ws.Cells.Item(223, 9).Value = "Primera"                              # I: Calidad
$ws.Cells.Item(223, 18).Value = "Hortaliza"                           # R: Clasificacion

# Match the date number format used by the rest of column D
$ws.Cells.Item(223, 4).NumberFormat = $ws.Cells.Item(222, 4).NumberFormat

